$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: 11/30/2023 (serial 45260), 3 hours, debugging visa requirements integration
$ws.Range("A66").Value = 45260
$ws.Range("A66").NumberFormat = "d-mmm"
$ws.Range("B66").Value = 3
$ws.Range("C66").Value = "trying to debug the visa requirements intergaration with the "

# Row 67: 12/2/2023 (serial 45262), 5 hours, weather API integrated, visa almost done
$ws.Range("A67").Value = 45262
$ws.Range("A67").NumberFormat = "d-mmm"
$ws.Range("B67").Value = 5
$ws.Range("C67").Value = "got the weather API intergrated into the UI, and have the VISA requirements almost working"

# Update the saved selection to match the new last cell, as Excel would after editing
$ws.Range("C67").Select() | Out-Null
